# Append two new log rows to the end of the log sheet (issue 43 fix):
# a "ReadyImage Signal Recieved" entry followed by a
# "ReadyPost Signal Recieved" entry, each stamped with its own timestamp,
# mirroring the existing ASP_SERVER log rows already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value = "2024 March 31 4:42:09 AM"
$ws.Cells.Item($newRow1, 2).Value = "ASP_SERVER"
$ws.Cells.Item($newRow1, 3).Value = "ReadyImage Signal Recieved"

$ws.Cells.Item($newRow2, 1).Value = "2024 March 31 4:42:10 AM"
$ws.Cells.Item($newRow2, 2).Value = "ASP_SERVER"
$ws.Cells.Item($newRow2, 3).Value = "ReadyPost Signal Recieved"
